$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D keeps its text semantics (values such as "1.00" or "0.999" must
# not be auto-coerced into numbers by Excel when assigned).
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "96.712.89"
$ws.Range("E2").Value = "  +0.29%  "

# Row 3
$ws.Range("D3").Value = "3.681.55"
$ws.Range("E3").Value = "  +0.44%  "

# Row 4
$ws.Range("E4").Value = "  -0.07%  "

# Row 5
$ws.Range("D5").Value = "235.96"
$ws.Range("E5").Value = "  -2.19%  "

# Row 6
$ws.Range("D6").Value = "1.93"
$ws.Range("E6").Value = "  +2.95%  "

# Row 7
$ws.Range("D7").Value = "653.18"
$ws.Range("E7").Value = "  -1.97%  "

# Row 8
$ws.Range("D8").Value = "0.422"
$ws.Range("E8").Value = "  -0.55%  "

# Row 9
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.06%  "

# Row 10
$ws.Range("D10").Value = "1.06"
$ws.Range("E10").Value = "  -1.76%  "

# Row 11
$ws.Range("D11").Value = "3.679.97"
$ws.Range("E11").Value = "  +0.47%  "

# Row 12
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "0.208"
$ws.Range("E12").Value = "  +1.90%  "

# Row 13
$ws.Range("B13").Value = "Avalanche"
$ws.Range("C13").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D13").Value = "43.91"
$ws.Range("E13").Value = "  -1.82%  "

# Row 14
$ws.Range("D14").Value = "0.0000296"
$ws.Range("E14").Value = "  +9.85%  "

# Row 15
$ws.Range("D15").Value = "6.72"
$ws.Range("E15").Value = "  +1.16%  "

# Row 16
$ws.Range("D16").Value = "4.368.72"
$ws.Range("E16").Value = "  +0.52%  "

# Row 17
$ws.Range("D17").Value = "96.435.77"
$ws.Range("E17").Value = "  +0.38%  "

# Row 18
$ws.Range("D18").Value = "8.92"
$ws.Range("E18").Value = "  +0.92%  "

# Row 19
$ws.Range("D19").Value = "3.664.04"
$ws.Range("E19").Value = "  +0.58%  "

# Row 20
$ws.Range("D20").Value = "12.95"
$ws.Range("E20").Value = "  +1.90%  "

# Row 21
$ws.Range("D21").Value = "18.54"
$ws.Range("E21").Value = "  +1.31%  "

# Row 22
$ws.Range("D22").Value = "0.508"
$ws.Range("E22").Value = "  -4.78%  "

# Row 23
$ws.Range("D23").Value = "519.56"
$ws.Range("E23").Value = "  -0.64%  "

# Row 24
$ws.Range("D24").Value = "3.39"
$ws.Range("E24").Value = "  -1.36%  "

# Row 25
$ws.Range("D25").Value = "0.0000210"
$ws.Range("E25").Value = "  +3.07%  "

# Row 26
$ws.Range("D26").Value = "6.87"
$ws.Range("E26").Value = "  -0.36%  "

# Row 27
$ws.Range("E27").Value = "  +22.12%  "

# Row 28
$ws.Range("D28").Value = "101.11"
$ws.Range("E28").Value = "  -1.33%  "

# Row 29
$ws.Range("D29").Value = "13.29"
$ws.Range("E29").Value = "  +2.76%  "

# Row 30
$ws.Range("D30").Value = "12.29"
$ws.Range("E30").Value = "  +1.67%  "

# Row 31
$ws.Range("D31").Value = "2.99"
$ws.Range("E31").Value = "  -1.41%  "

# Row 33
$ws.Range("E33").Value = "  +1.59%  "

# Row 34
$ws.Range("D34").Value = "1.83"
$ws.Range("E34").Value = "  +1.10%  "

# Row 35
$ws.Range("E35").Value = "  +0.04%  "

# Row 36
$ws.Range("D36").Value = "32.05"
$ws.Range("E36").Value = "  -1.88%  "

# Row 37
$ws.Range("D37").Value = "641.39"
$ws.Range("E37").Value = "  +3.92%  "

# Row 38
$ws.Range("D38").Value = "0.589"
$ws.Range("E38").Value = "  +0.90%  "

# Row 39
$ws.Range("D39").Value = "8.73"
$ws.Range("E39").Value = "  +0.13%  "

# Row 40
$ws.Range("E40").Value = "  +0.00%  "

# Row 41
$ws.Range("D41").Value = "6.78"
$ws.Range("E41").Value = "  +10.45%  "

# Row 42
$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").Value = "40.77"
$ws.Range("E42").Value = "  -4.71%  "

# Row 43
$ws.Range("B43").Value = "ImmutableX"
$ws.Range("C43").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D43").Value = "2.03"
$ws.Range("E43").Value = "  +5.06%  "

# Row 44
$ws.Range("E44").Value = "  +0.78%  "

# Row 45
$ws.Range("D45").Value = "0.948"
$ws.Range("E45").Value = "  -0.14%  "

# Row 46
$ws.Range("D46").Value = "0.453"
$ws.Range("E46").Value = "  +6.23%  "

# Row 47
$ws.Range("D47").Value = "0.0454"
$ws.Range("E47").Value = "  +0.12%  "

# Row 48
$ws.Range("D48").Value = "23.61"

# Row 49
$ws.Range("D49").Value = "2.26"
$ws.Range("E49").Value = "  -0.97%  "

# Row 50
$ws.Range("D50").Value = "8.47"
$ws.Range("E50").Value = "  +0.09%  "

# Row 51
$ws.Range("D51").Value = "3.53"
$ws.Range("E51").Value = "  -0.48%  "
